$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C4").Value = -12.974
$ws.Range("D4").Value = -8.25

$ws.Range("D5").Value = -8.581999999999999

$ws.Range("C6").Value = -12.632
$ws.Range("D6").Value = -8.31

$ws.Range("C7").Value = -13.214

$ws.Range("C8").Value = -12.694
$ws.Range("D8").Value = -8.385999999999999

$ws.Range("C16").Value = -12.715
$ws.Range("D16").Value = -8.568000000000001

$ws.Range("C20").Value = -12.9

$ws.Range("C21").Value = -13.214

$ws.Range("D22").Value = -8.177000000000001
